# Add 2022-Q1 data: insert a new "2022-Q1" sheet with per-fund holdings before the
# "总计" (totals) summary sheet, and add a corresponding summary row to "总计".

$excel.DisplayAlerts = $false
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write a value as TEXT (preserving exact formatting / leading zeros),
# matching the source workbook's convention of storing numeric-looking figures
# (fund codes, percentages, AUM, etc.) as plain strings rather than numbers.
# ---------------------------------------------------------------------------
function Set-TextValue($range, [string]$text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# ---------------------------------------------------------------------------
# Step 1: remove the existing "总计" sheet. It will be recreated at the end of
# the workbook so that it receives a fresh (higher) sheetId, matching the new
# ordering: ... 2021-Q4, 2022-Q1, 总计.
# ---------------------------------------------------------------------------
$oldTotal = $wb.Worksheets.Item("总计")
$oldTotal.Delete()

# ---------------------------------------------------------------------------
# Step 2: create the new "2022-Q1" sheet right after "2021-Q4".
# ---------------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$q1New = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $q4)
$q1New.Name = "2022-Q1"

# ---------------------------------------------------------------------------
# Step 3: re-create the "总计" sheet right after "2022-Q1".
# ---------------------------------------------------------------------------
$totalNew = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $q1New)
$totalNew.Name = "总计"

# ---------------------------------------------------------------------------
# Step 4: populate "2022-Q1" with per-fund holding data, copying the header /
# index-column style (bold, bordered) from the "2021-Q4" sheet which already
# uses the same layout.
# ---------------------------------------------------------------------------
Set-TextValue $q1New.Range("B1") "基金代码"
Set-TextValue $q1New.Range("C1") "基金名称"
Set-TextValue $q1New.Range("D1") "基金规模"
Set-TextValue $q1New.Range("E1") "股票总仓位"
Set-TextValue $q1New.Range("F1") "仓位占比"
Set-TextValue $q1New.Range("G1") "持有市值(亿元)"
Set-TextValue $q1New.Range("H1") "仓位排名"

$q4.Range("B1:H1").Copy()
$q1New.Range("B1:H1").PasteSpecial(-4122) | Out-Null

$fundRows = @(
    @{ A = 0; B = "007751"; C = "景顺长城中证沪港深红利成长低波动指数A"; D = "0.83"; E = "91.29"; F = "3.45"; G = "0.0286"; H = 1 },
    @{ A = 1; B = "005702"; C = "恒生前海港股通高股息低波动指数";         D = "0.29"; E = "94.14"; F = "2.32"; G = "0.0067"; H = 9 },
    @{ A = 2; B = "007760"; C = "景顺长城中证沪港深红利成长低波动指数C"; D = "0.06"; E = "91.29"; F = "3.45"; G = "0.0021"; H = 1 }
)

$r = 2
foreach ($row in $fundRows) {
    $q1New.Cells.Item($r, 1).Value = $row.A
    $q4.Range("A2").Copy()
    $q1New.Cells.Item($r, 1).PasteSpecial(-4122) | Out-Null
    $q1New.Cells.Item($r, 1).Value = $row.A

    Set-TextValue $q1New.Cells.Item($r, 2) $row.B
    Set-TextValue $q1New.Cells.Item($r, 3) $row.C
    Set-TextValue $q1New.Cells.Item($r, 4) $row.D
    Set-TextValue $q1New.Cells.Item($r, 5) $row.E
    Set-TextValue $q1New.Cells.Item($r, 6) $row.F
    Set-TextValue $q1New.Cells.Item($r, 7) $row.G
    $q1New.Cells.Item($r, 8).Value = $row.H

    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Step 5: populate "总计" with the historical summary rows plus the new
# 2022-Q1 entry inserted at the top.
# ---------------------------------------------------------------------------
Set-TextValue $totalNew.Range("B1") "日期"
Set-TextValue $totalNew.Range("C1") "持有数量(只)"
Set-TextValue $totalNew.Range("D1") "持有市值(亿元)"

$q4.Range("B1:D1").Copy()
$totalNew.Range("B1:D1").PasteSpecial(-4122) | Out-Null

$totalRows = @(
    @{ A = 0; B = "2022-Q1"; C = 3; D = 0.04 },
    @{ A = 1; B = "2021-Q4"; C = 4; D = 2.36 },
    @{ A = 2; B = "2021-Q3"; C = 1; D = 2.55 },
    @{ A = 3; B = "2021-Q2"; C = 3; D = 1.76 },
    @{ A = 4; B = "2021-Q1"; C = 3; D = 1.51 },
    @{ A = 5; B = "2020-Q4"; C = 6; D = 0.26 }
)

$r = 2
foreach ($row in $totalRows) {
    $totalNew.Cells.Item($r, 1).Value = $row.A
    $q4.Range("A2").Copy()
    $totalNew.Cells.Item($r, 1).PasteSpecial(-4122) | Out-Null
    $totalNew.Cells.Item($r, 1).Value = $row.A

    Set-TextValue $totalNew.Cells.Item($r, 2) $row.B
    $totalNew.Cells.Item($r, 3).Value = $row.C
    $totalNew.Cells.Item($r, 4).Value = $row.D

    $r = $r + 1
}
